{"js": "// Replace each three-digit-by-one-digit multiplication problem's\n// operands with the updated values from the commit, matching each\n// old expression exactly once (all problem strings are unique in the\n// document) and preserving the surrounding run formatting.\nconst replacements = [\n  [\"846\u00d74=\", \"788\u00d73=\"],\n  [\"172\u00d79=\", \"308\u00d75=\"],\n  [\"398\u00d74=\", \"634\u00d75=\"],\n  [\"419\u00d77=\", \"170\u00d73=\"],\n  [\"613\u00d77=\", \"137\u00d76=\"],\n  [\"965\u00d74=\", \"831\u00d72=\"],\n  [\"393\u00d75=\", \"122\u00d73=\"],\n  [\"179\u00d79=\", \"304\u00d73=\"],\n  [\"941\u00d76=\", \"581\u00d74=\"],\n  [\"243\u00d72=\", \"991\u00d76=\"],\n  [\"494\u00d74=\", \"928\u00d72=\"],\n  [\"750\u00d78=\", \"346\u00d77=\"],\n  [\"749\u00d74=\", \"402\u00d75=\"],\n  [\"621\u00d75=\", \"259\u00d73=\"],\n  [\"750\u00d72=\", \"378\u00d77=\"],\n  [\"416\u00d72=\", \"840\u00d73=\"],\n  [\"242\u00d74=\", \"205\u00d75=\"],\n  [\"580\u00d72=\", \"230\u00d77=\"],\n  [\"584\u00d72=\", \"299\u00d77=\"],\n  [\"782\u00d76=\", \"380\u00d74=\"],\n  [\"816\u00d74=\", \"203\u00d74=\"],\n  [\"131\u00d72=\", \"968\u00d74=\"],\n  [\"748\u00d76=\", \"798\u00d79=\"],\n  [\"709\u00d74=\", \"901\u00d73=\"],\n  [\"787\u00d77=\", \"958\u00d78=\"]\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Could not find text to replace: ${oldText}`);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update each three-digit-by-one-digit multiplication problem's\n# operands to the new values from the commit. Every \"old\" expression\n# occurs exactly once in the document, so a ReplaceAll Find/Replace\n# for each pair maps 1:1 onto the diff's per-cell text changes while\n# leaving all other run/paragraph formatting untouched.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"846\u00d74=\"; New = \"788\u00d73=\" },\n    @{ Old = \"172\u00d79=\"; New = \"308\u00d75=\" },\n    @{ Old = \"398\u00d74=\"; New = \"634\u00d75=\" },\n    @{ Old = \"419\u00d77=\"; New = \"170\u00d73=\" },\n    @{ Old = \"613\u00d77=\"; New = \"137\u00d76=\" },\n    @{ Old = \"965\u00d74=\"; New = \"831\u00d72=\" },\n    @{ Old = \"393\u00d75=\"; New = \"122\u00d73=\" },\n    @{ Old = \"179\u00d79=\"; New = \"304\u00d73=\" },\n    @{ Old = \"941\u00d76=\"; New = \"581\u00d74=\" },\n    @{ Old = \"243\u00d72=\"; New = \"991\u00d76=\" },\n    @{ Old = \"494\u00d74=\"; New = \"928\u00d72=\" },\n    @{ Old = \"750\u00d78=\"; New = \"346\u00d77=\" },\n    @{ Old = \"749\u00d74=\"; New = \"402\u00d75=\" },\n    @{ Old = \"621\u00d75=\"; New = \"259\u00d73=\" },\n    @{ Old = \"750\u00d72=\"; New = \"378\u00d77=\" },\n    @{ Old = \"416\u00d72=\"; New = \"840\u00d73=\" },\n    @{ Old = \"242\u00d74=\"; New = \"205\u00d75=\" },\n    @{ Old = \"580\u00d72=\"; New = \"230\u00d77=\" },\n    @{ Old = \"584\u00d72=\"; New = \"299\u00d77=\" },\n    @{ Old = \"782\u00d76=\"; New = \"380\u00d74=\" },\n    @{ Old = \"816\u00d74=\"; New = \"203\u00d74=\" },\n    @{ Old = \"131\u00d72=\"; New = \"968\u00d74=\" },\n    @{ Old = \"748\u00d76=\"; New = \"798\u00d79=\" },\n    @{ Old = \"709\u00d74=\"; New = \"901\u00d73=\" },\n    @{ Old = \"787\u00d77=\"; New = \"958\u00d78=\" }\n)\n\n$wdReplaceAll = 2\n$wdFindContinue = 1\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair.Old\n    $find.Replacement.Text = $pair.New\n    $find.Execute($pair.Old, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $pair.New, $wdReplaceAll) | Out-Null\n}\n"}
